# Add files via upload
# - Insert a new player row (M. Kowalczyk) before J. Kuzma (old row 19)
# - Refresh season-stat totals for every player after a new match
# - Narrow column F and update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new player by inserting a row at position 19.
$ws.Rows.Item(19).Insert()

# 2) Full target data set (rows 2..35, columns A..M) after the roster/stat update.

$data = @(
    ,@('D. Arndt', 20, 'B', 270, 3, 3, 0, 0, 12, 0, 0, 0, 0)
    ,@('A. Bobek', 17, 'B', 0, 0, 0, 0, 0, 2, 0, 0, 0, 0)
    ,@('M. Kot', 19, 'B', 90, 1, 1, 0, 0, 5, 0, 0, 0, 0)
    ,@('M. Kozioł', 33, 'B', 1620, 18, 18, 0, 0, 1, 0, 1, 0, 0)
    ,@('M. Kołba', 29, 'B', 0, 0, 0, 0, 0, 4, 0, 0, 0, 0)
    ,@('M. Bąkowicz', 20, 'O', 1150, 15, 14, 1, 4, 5, 0, 2, 0, 0)
    ,@('M. Dąbrowski', 34, 'O', 830, 11, 11, 0, 3, 1, 1, 3, 1, 0)
    ,@('K. Dankowski', 25, 'O', 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    ,@('A. Klimczak', 24, 'O', 1011, 14, 10, 4, 1, 8, 0, 2, 0, 0)
    ,@('O. Koprowski', 22, 'O', 1050, 14, 11, 3, 1, 8, 0, 4, 0, 0)
    ,@('M. Lorenc', 20, 'O', 377, 5, 4, 1, 0, 6, 0, 1, 0, 0)
    ,@('A. Marciniak', 33, 'O', 1019, 12, 12, 0, 2, 5, 0, 2, 0, 0)
    ,@('Nacho Monsalve', 27, 'O', 257, 6, 3, 3, 1, 5, 2, 2, 0, 0)
    ,@('B. Szeliga', 29, 'O', 1394, 18, 15, 3, 2, 3, 1, 4, 0, 0)
    ,@('M. Wolski', 24, 'O', 1080, 20, 14, 6, 13, 7, 2, 4, 0, 0)
    ,@('Antonio Domínguez', 28, 'P', 1062, 14, 12, 2, 5, 2, 2, 2, 0, 0)
    ,@('P. Gryszkiewicz', 20, 'P', 236, 10, 1, 9, 1, 13, 0, 2, 0, 0)
    ,@('M. Kowalczyk', 17, 'P', 50, 2, 0, 2, 0, 2, 0, 0, 0, 0)
    ,@('J. Kuźma', 18, 'P', 221, 5, 2, 3, 1, 9, 0, 0, 0, 0)
    ,@('M. Lipien', 18, 'P', 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    ,@('Javi Moreno', 24, 'P', 1210, 18, 13, 5, 5, 5, 0, 2, 0, 0)
    ,@('D. Nowacki', 23, 'P', 39, 4, 0, 4, 0, 6, 0, 0, 0, 0)
    ,@('M. Rozwandowicz', 27, 'P', 1036, 16, 13, 3, 11, 6, 1, 5, 0, 0)
    ,@('Pirulo', 29, 'P', 1580, 20, 19, 1, 7, 1, 5, 5, 0, 0)
    ,@('J. Tosik', 34, 'P', 607, 15, 6, 9, 6, 10, 0, 4, 0, 0)
    ,@('M. Trąbka', 24, 'P', 1071, 13, 13, 0, 4, 0, 0, 1, 0, 0)
    ,@('M. Wszołek', 19, 'P', 0, 0, 0, 0, 0, 6, 0, 0, 0, 0)
    ,@('Ricardinho', 32, 'N', 701, 12, 7, 5, 6, 6, 0, 1, 0, 1)
    ,@('Samu Corral', 29, 'N', 376, 10, 3, 7, 2, 9, 0, 1, 0, 0)
    ,@('D. Gmosiński', 21, 'N', 0, 0, 0, 0, 0, 2, 0, 0, 0, 0)
    ,@('K. Ibe-Torti', 20, 'N', 306, 9, 3, 6, 2, 12, 0, 1, 0, 0)
    ,@('P. Janczukowicz', 22, 'N', 407, 11, 5, 6, 5, 6, 3, 3, 0, 0)
    ,@('S. Jurić', 23, 'N', 619, 14, 7, 7, 6, 9, 2, 2, 0, 0)
    ,@('M. Radaszkiewicz', 24, 'N', 740, 11, 8, 3, 7, 3, 3, 0, 0, 0)
)

# 3) Write every row/column back out.
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $rowNum = $i + 2
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}

# 4) Column F narrowed after the stat refresh - autofit to the new content.
$ws.Columns.Item(6).AutoFit()

# 5) Restore the active window state recorded after the edit.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("A38:XFD72").Select()
